$p = $ppt.ActivePresentation

# 1) Slide 6: fix wording "Bottom chart" -> "Top chart" in the bullet list shape,
#    replacing the full run text so the run stays intact (no split).
$s6 = $p.Slides.Item(6)
$shape6 = $s6.Shapes.Item(7)
$tr = $shape6.TextFrame.TextRange
$oldText = "Bottom chart represents the Average Overall Happiness to Freedom factor scores for 2022."
$newText = "Top chart represents the Average Overall Happiness to Freedom factor scores for 2022."
$run = $tr.Characters(1, $oldText.Length)
$run.Text = $newText

# 2) Slide 4: nudge picture "Picture 8" horizontally.
#    Shape.Left/.Top are in points; the OOXML offsets are in EMUs (1 pt = 12700 EMU).
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(7)
$shape4.Left = 7728087 / 12700
